$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 22089.56346396808
$ws.Range("D2").Value = 0
$ws.Range("G2").Value = -2332529839.07643

# Row 3
$ws.Range("B3").Value = 22444.79953830915
$ws.Range("D3").Value = 0
$ws.Range("G3").Value = -3104218115.362121

# Row 4
$ws.Range("B4").Value = 23786.41279629565
$ws.Range("D4").Value = 0
$ws.Range("G4").Value = -2373274113.218047

# Row 5
$ws.Range("B5").Value = 24161.1029046647
$ws.Range("D5").Value = 0
$ws.Range("G5").Value = -3180384352.007998

# Row 6
$ws.Range("B6").Value = 24066.91794378401
$ws.Range("D6").Value = 0
$ws.Range("G6").Value = -29762704211.95647

# Row 7
$ws.Range("B7").Value = 24357.23370853326
$ws.Range("D7").Value = 0
$ws.Range("G7").Value = -30640596520.94985

# Row 8
$ws.Range("B8").Value = 23440.34122279853
$ws.Range("D8").Value = 0
$ws.Range("G8").Value = -29006713587.66103

# Row 9
$ws.Range("B9").Value = 23724.10495496288
$ws.Range("D9").Value = 0
$ws.Range("G9").Value = -29863687044.67785
